# Add a list of functions used in solving problems:
# 문자열.strip() / 문자열.lstrip() / 문자열.rstrip()

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: strip()
$ws.Range("A20").Value = "문자열.strip()"
$ws.Range("B20").Value = "문자열의 양쪽 끝에서 특정 문자를 제거"
$ws.Range("C20").Value = "기본적으로 공백(스페이스, 탭, 개행 등)을 제거하지만, 제거할 문자(또는 문자열)를 지정 가능"
$ws.Range("D20").Value = "string.strip([chars])"
$ws.Rows.Item(20).RowHeight = 40

# Row 21: lstrip()
$ws.Range("A21").Value = "문자열.lstrip()"
$ws.Range("B21").Value = "문자열의 왼쪽 끝에서만 지정한 문자 제거"
$ws.Range("B21").Font.Name = "D2Coding"
$ws.Range("B21").Font.Size = 12
$ws.Range("B21").Font.Color = 921102
$ws.Rows.Item(21).RowHeight = 40

# Row 22: rstrip()
$ws.Range("A22").Value = "문자열.rstrip()"
$ws.Range("B22").Value = "문자열의 오른쪽 끝에서만 지정한 문자 제거"
$ws.Rows.Item(22).RowHeight = 40

# Update selection to match final cursor position
$ws.Range("B22").Select()
